$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in the Price column hold text-formatted numbers (e.g. "26.277.23",
# "217.56") rather than numeric values, so each is switched to Text format
# before the value is written to prevent Excel from re-interpreting the
# string as a number.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.277.23'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.688.35'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.56'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5376'
$ws.Range('E6').Value = '  +2.28%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2720'
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06418'
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.62'
$ws.Range('E10').Value = '  -1.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07665'
$ws.Range('E11').Value = '  +1.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.687.73'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.523'
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5775'
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008372'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.58'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.306.42'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.008'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.881'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.83'
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '190.47'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.252'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.59'
$ws.Range('E24').Value = '  +2.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1293'
$ws.Range('E25').Value = '  +3.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.849'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.80'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06204'
$ws.Range('E28').Value = '  -4.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.377'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.598'
$ws.Range('E31').Value = '  -0.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.578'
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.675'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.027'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6184'
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.423'
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.758'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01649'
$ws.Range('E38').Value = '  +1.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.122'
$ws.Range('E39').Value = '  -4.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.107.07'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8802'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.012'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.94'
$ws.Range('E43').Value = '  +0.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.839.62'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000112'
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05282'
$ws.Range('E49').Value = '  +0.24%  '
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.057'
$ws.Range('E51').Value = '  -0.44%  '
